$d = $word.ActiveDocument

# --- Change 1 ---------------------------------------------------------
# Remove the stray space-only run that precedes "jmbgFrom" and make the
# "jmbgFrom" run use hr-HR (Croatian) instead of en-US, keeping its other
# direct formatting (Times New Roman / 14pt) untouched.
$rng = $d.Content
$found = $rng.Find.Execute("jmbgFrom", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    # The run immediately before "jmbgFrom" is a single space with only
    # w:lang in its rPr - delete it entirely (text + run).
    $spaceRange = $d.Range($rng.Start - 1, $rng.Start)
    if ($spaceRange.Text -eq " ") {
        $spaceRange.Delete()
    }

    # Re-find "jmbgFrom" now that the space is gone, then retarget its
    # language to hr-HR.
    $rng2 = $d.Content
    $found2 = $rng2.Find.Execute("jmbgFrom", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2) {
        $rng2.LanguageID = "hr-HR"
    }
}

# --- Change 2 ---------------------------------------------------------
# Insert a new run containing "Num" right after the existing "jmbg" run
# (the one inside "JMBG: {jmbg}"), using the same direct formatting
# (Times New Roman / 14pt / en-US) as that run, but kept as its own
# separate <w:r> rather than merged into "jmbg".
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("{jmbg}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $insertPos = $rng3.End - 1
    $insertPoint = $d.Range($insertPos, $insertPos)
    $insertPoint.InsertAfter("Num")

    $numRange = $d.Range($insertPos, $insertPos + 3)
    # Toggle a property on/off so the new text keeps its own run instead
    # of silently being coalesced with the preceding identically
    # formatted "jmbg" run.
    $numRange.Bold = 1
    $numRange.Bold = 0
}
